$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan de Pruebas")

# Row 6 (Test Case 5) "Estado" column F went from "PTE" to "OK".
# Copy the formatting used by the other "OK" cells (e.g. F5) onto F6,
# then update its value/text to "OK".
$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("F6").Value = "OK"

$excel.CutCopyMode = 0
